# Apply updated monitoring data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: total_registros for INCIO SANCHEZ PAOLA KATHERINE increases 89 -> 90
$ws.Range("B2").Value = 90

# Rows 7-15 are reordered (by descending total_registros) and/or have
# updated counts. Rewrite the A/B values for rows 7 through 15 to match
# the new ranking.
$ws.Range("A7").Value = "DELGADO VASQUEZ FLOR MAGALY"
$ws.Range("B7").Value = 77

$ws.Range("A8").Value = "PEREZ LINARES TATHIANA"
$ws.Range("B8").Value = 77

$ws.Range("A9").Value = "MEDINA TAPIA ANA YULI"
$ws.Range("B9").Value = 76

$ws.Range("A10").Value = "MONDRAGON HERNANDEZ WILMER JUNIOR"
$ws.Range("B10").Value = 76

$ws.Range("A11").Value = "CAMPOS PEREZ YOVERLY"
$ws.Range("B11").Value = 75

$ws.Range("A12").Value = "CHAVEZ VILLANUEVA SILVIA JANETH"
$ws.Range("B12").Value = 74

$ws.Range("A13").Value = "LOZADA ROJAS LUZ ELENA"
$ws.Range("B13").Value = 73

$ws.Range("A14").Value = "VASQUEZ SILVA ALOIS ADOLF"
$ws.Range("B14").Value = 73

$ws.Range("A15").Value = "SOTO LOZANO LUZDINA"
$ws.Range("B15").Value = 70
